$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The handback took place and is now in sync with en-US: update the status text
# everywhere it appears (Overview rollup + each language sheet).
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

# --- zh-cn sheet: record the handback target/handback files + datetime ---
$wsZhCn.Range("E2").Value = "658f4757-a488-4c31-8103-b96abf108c5b.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a2509e2722a9d634a9921797a45f8e86909878fb/e2e/658f4757-a488-4c31-8103-b96abf108c5b.md", "", "", "658f4757-a488-4c31-8103-b96abf108c5b.md")
$wsZhCn.Range("E2").Style = "HyperLink"

$wsZhCn.Range("F2").Value = "658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ff0e93870a34d416c3a57ded7cd64a125adc5e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.zh-cn.xlf", "", "", "658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.zh-cn.xlf")
$wsZhCn.Range("F2").Style = "HyperLink"

$wsZhCn.Range("G2").Value = "2016-03-02 07:07:12"

$wsZhCn.Range("E3").Value = "fe1add01-6e6d-4054-82d4-5637c1f53052.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a2509e2722a9d634a9921797a45f8e86909878fb/e2e/fe1add01-6e6d-4054-82d4-5637c1f53052.md", "", "", "fe1add01-6e6d-4054-82d4-5637c1f53052.md")
$wsZhCn.Range("E3").Style = "HyperLink"

$wsZhCn.Range("F3").Value = "fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ff0e93870a34d416c3a57ded7cd64a125adc5e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.zh-cn.xlf", "", "", "fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.zh-cn.xlf")
$wsZhCn.Range("F3").Style = "HyperLink"

$wsZhCn.Range("G3").Value = "2016-03-02 07:07:12"

# --- de-de sheet: same handback bookkeeping ---
$wsDeDe.Range("E2").Value = "658f4757-a488-4c31-8103-b96abf108c5b.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a2509e2722a9d634a9921797a45f8e86909878fb/e2e/658f4757-a488-4c31-8103-b96abf108c5b.md", "", "", "658f4757-a488-4c31-8103-b96abf108c5b.md")
$wsDeDe.Range("E2").Style = "HyperLink"

$wsDeDe.Range("F2").Value = "658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/235ea150bf8d216ed8f4993807b9deeb6657694c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.de-de.xlf", "", "", "658f4757-a488-4c31-8103-b96abf108c5b.59ed6bb3e0b1f45bcb834b274f35b8c8674bdb53.de-de.xlf")
$wsDeDe.Range("F2").Style = "HyperLink"

$wsDeDe.Range("G2").Value = "2016-03-02 07:07:31"

$wsDeDe.Range("E3").Value = "fe1add01-6e6d-4054-82d4-5637c1f53052.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a2509e2722a9d634a9921797a45f8e86909878fb/e2e/fe1add01-6e6d-4054-82d4-5637c1f53052.md", "", "", "fe1add01-6e6d-4054-82d4-5637c1f53052.md")
$wsDeDe.Range("E3").Style = "HyperLink"

$wsDeDe.Range("F3").Value = "fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/235ea150bf8d216ed8f4993807b9deeb6657694c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.de-de.xlf", "", "", "fe1add01-6e6d-4054-82d4-5637c1f53052.70f9ac5d14508eef1fe499bac3b56eb24e3d0570.de-de.xlf")
$wsDeDe.Range("F3").Style = "HyperLink"

$wsDeDe.Range("G3").Value = "2016-03-02 07:07:31"
